$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$c = $ws.Columns("E:G")
Write-Host $c
